$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.645.44'
$ws.Range('E2').Value = '  +4.34%  '
$ws.Range('D3').Value = '1.870.37'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.45'
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4700'
$ws.Range('E7').Value = '  +3.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3994'
$ws.Range('E8').Value = '  +5.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.67'
$ws.Range('E9').Value = '  +3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08044'
$ws.Range('E10').Value = '  +2.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.002'
$ws.Range('E11').Value = '  +3.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.11'
$ws.Range('E12').Value = '  +5.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.054'
$ws.Range('E13').Value = '  +3.34%  '
$ws.Range('D14').Value = '1.862.03'
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('E15').Value = '  +3.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.30'
$ws.Range('E16').Value = '  +2.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001044'
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06607'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.62'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.51%  '
$ws.Range('D22').Value = '28.635.89'
$ws.Range('E22').Value = '  +4.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.507'
$ws.Range('E23').Value = '  +3.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.06'
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.255'
$ws.Range('E25').Value = '  -2.05%  '
$ws.Range('D26').Value = '2.086.71'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.69'
$ws.Range('E27').Value = '  +2.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.74'
$ws.Range('E28').Value = '  +1.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.125'
$ws.Range('E29').Value = '  +3.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.485'
$ws.Range('E30').Value = '  +4.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.61'
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9784'
$ws.Range('E32').Value = '  +3.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09536'
$ws.Range('E33').Value = '  +2.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.673'
$ws.Range('E34').Value = '  +2.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.385'
$ws.Range('E35').Value = '  +5.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.370'
$ws.Range('E36').Value = '  +2.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06212'
$ws.Range('E37').Value = '  +5.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02255'
$ws.Range('E38').Value = '  +3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.457'
$ws.Range('E39').Value = '  +5.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.180'
$ws.Range('E40').Value = '  +2.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5949'
$ws.Range('E41').Value = '  +3.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9988'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1886'
$ws.Range('E43').Value = '  +3.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.33'
$ws.Range('E44').Value = '  +3.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.259'
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5573'
$ws.Range('E46').Value = '  +2.44%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07439'
$ws.Range('E47').Value = '  +12.62%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.12'
$ws.Range('E48').Value = '  +1.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.960'
$ws.Range('E49').Value = '  +5.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.076'
$ws.Range('E50').Value = '  +13.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '112.42'
$ws.Range('E51').Value = '  +1.98%  '
